# Added average stats after subsampling
#
# This script updates the "Metrics Table" worksheet:
#   - Header labels "Training Metrics" / "Test Metrics" become
#     "Average Training Metrics" / "Average Test Metrics" (values that were
#     presumably computed over a larger run are replaced by the average
#     over several subsamples).
#   - The single-sample "Linear" kernel row (row 3) no longer carries its
#     own Accuracy figures - it is blanked out.
#   - All the numeric metrics throughout the table are refreshed with the
#     new averaged values.
#   - Two previously-empty rows ("Degree 15" / "Degree 20") now get their
#     own averaged metrics filled in.
#   - The active-cell selection moves from A10 to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: rename the merged-header labels ------------------------
$ws.Range("B1").Value = "Average Training Metrics"
$ws.Range("F1").Value = "Average Test Metrics"

# --- Row 3 ("Linear"): drop the old single-run figures ----------------------
$ws.Range("B3").Clear()
$ws.Range("F3").Clear()

# --- Row 4 ("RBF") -----------------------------------------------------------
$ws.Range("B4").Value = 0.73
$ws.Range("C4").Value = 0.74
$ws.Range("D4").Value = 0.84
$ws.Range("E4").Value = 0.78
$ws.Range("F4").Value = 0.6
$ws.Range("G4").Value = 0.6
$ws.Range("H4").Value = 0.68
$ws.Range("I4").Value = 0.6

# --- Row 6 ("Degree 2") ------------------------------------------------------
$ws.Range("B6").Value = 0.64
$ws.Range("C6").Value = 0.67
$ws.Range("D6").Value = 0.85
$ws.Range("E6").Value = 0.72
$ws.Range("F6").Value = 0.53
$ws.Range("G6").Value = 0.53
$ws.Range("H6").Value = 0.62
$ws.Range("I6").Value = 0.55000000000000004

# --- Row 7 ("Degree 3") ------------------------------------------------------
$ws.Range("B7").Value = 0.67
$ws.Range("C7").Value = 0.66
$ws.Range("D7").Value = 0.9
$ws.Range("E7").Value = 0.76
$ws.Range("F7").Value = 0.59
$ws.Range("G7").Value = 0.57999999999999996
$ws.Range("H7").Value = 0.85
$ws.Range("I7").Value = 0.69

# --- Row 8 ("Degree 4") ------------------------------------------------------
$ws.Range("B8").Value = 0.66
$ws.Range("C8").Value = 0.7
$ws.Range("D8").Value = 0.82
$ws.Range("E8").Value = 0.7
$ws.Range("F8").Value = 0.56999999999999995
$ws.Range("G8").Value = 0.57999999999999996
$ws.Range("H8").Value = 0.65
$ws.Range("I8").Value = 0.56999999999999995

# --- Row 9 ("Degree 5") ------------------------------------------------------
$ws.Range("B9").Value = 0.68
$ws.Range("C9").Value = 0.72
$ws.Range("D9").Value = 0.83
$ws.Range("E9").Value = 0.73
$ws.Range("F9").Value = 0.59
$ws.Range("G9").Value = 0.56999999999999995
$ws.Range("H9").Value = 0.73
$ws.Range("I9").Value = 0.61

# --- Row 10 ("Degree 6") -----------------------------------------------------
$ws.Range("B10").Value = 0.68
$ws.Range("C10").Value = 0.73
$ws.Range("D10").Value = 0.82
$ws.Range("E10").Value = 0.72
$ws.Range("F10").Value = 0.56000000000000005
$ws.Range("G10").Value = 0.53
$ws.Range("H10").Value = 0.62
$ws.Range("I10").Value = 0.54

# --- Row 11 ("Degree 7") -----------------------------------------------------
$ws.Range("B11").Value = 0.7
$ws.Range("C11").Value = 0.74
$ws.Range("D11").Value = 0.84
$ws.Range("E11").Value = 0.74
$ws.Range("F11").Value = 0.61
$ws.Range("G11").Value = 0.56999999999999995
$ws.Range("H11").Value = 0.73
$ws.Range("I11").Value = 0.6

# --- Row 12 ("Degree 8") -----------------------------------------------------
$ws.Range("B12").Value = 0.7
$ws.Range("C12").Value = 0.74
$ws.Range("D12").Value = 0.84
$ws.Range("E12").Value = 0.74
$ws.Range("F12").Value = 0.56999999999999995
$ws.Range("G12").Value = 0.59
$ws.Range("H12").Value = 0.6
$ws.Range("I12").Value = 0.55000000000000004

# --- Row 13 ("Degree 9") -----------------------------------------------------
$ws.Range("B13").Value = 0.71
$ws.Range("C13").Value = 0.76
$ws.Range("D13").Value = 0.84
$ws.Range("E13").Value = 0.75
$ws.Range("F13").Value = 0.57999999999999996
$ws.Range("G13").Value = 0.5
$ws.Range("H13").Value = 0.69
$ws.Range("I13").Value = 0.56000000000000005

# --- Row 14 ("Degree 10") ----------------------------------------------------
$ws.Range("B14").Value = 0.71
$ws.Range("C14").Value = 0.77
$ws.Range("D14").Value = 0.85
$ws.Range("E14").Value = 0.75
$ws.Range("F14").Value = 0.56999999999999995
$ws.Range("G14").Value = 0.56999999999999995
$ws.Range("H14").Value = 0.6
$ws.Range("I14").Value = 0.55000000000000004

# --- Row 15 ("Degree 15") - previously blank, now filled in ------------------
$ws.Range("B15").Value = 0.73
$ws.Range("C15").Value = 0.78
$ws.Range("D15").Value = 0.85
$ws.Range("E15").Value = 0.77
$ws.Range("F15").Value = 0.57999999999999996
$ws.Range("G15").Value = 0.53
$ws.Range("H15").Value = 0.67
$ws.Range("I15").Value = 0.56000000000000005

# --- Row 16 ("Degree 20") - previously blank, now filled in ------------------
$ws.Range("B16").Value = 0.73
$ws.Range("C16").Value = 0.78
$ws.Range("D16").Value = 0.85
$ws.Range("E16").Value = 0.77
$ws.Range("F16").Value = 0.55000000000000004
$ws.Range("G16").Value = 0.54
$ws.Range("H16").Value = 0.6
$ws.Range("I16").Value = 0.54

# --- Rows 18-23 ("Custom" / p = 1..10): only the Test Metrics shift slightly -
$ws.Range("F18").Value = 0.56999999999999995
$ws.Range("G18").Value = 0.56999999999999995
$ws.Range("I18").Value = 0.72

$ws.Range("F19").Value = 0.56999999999999995
$ws.Range("G19").Value = 0.56999999999999995
$ws.Range("I19").Value = 0.72

$ws.Range("F20").Value = 0.56999999999999995
$ws.Range("G20").Value = 0.56999999999999995
$ws.Range("I20").Value = 0.72

$ws.Range("F21").Value = 0.56999999999999995
$ws.Range("G21").Value = 0.56999999999999995
$ws.Range("I21").Value = 0.72

$ws.Range("F22").Value = 0.56999999999999995
$ws.Range("G22").Value = 0.56999999999999995
$ws.Range("I22").Value = 0.72

$ws.Range("F23").Value = 0.56999999999999995
$ws.Range("G23").Value = 0.56999999999999995
$ws.Range("I23").Value = 0.72

# --- Move the active selection from A10 to A8 --------------------------------
$ws.Range("A8").Select()
